# chore: add user name in checkout
#
# - Append a new checkout timestamp to the existing check_in value in C4.
# - Add a new visitor row (row 5) for "Aahil Alwani".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: append the new check-in timestamp onto the existing one.
$ws.Range("C4").Value = "2025-05-02T11:32:00.422Z | 2025-05-05T08:08:42.144672Z"

# Row 5: new visitor record.
# Format column B as text first so the long, numeric-looking CNIC value is
# stored verbatim instead of being coerced into a floating point number.
$ws.Range("B5").NumberFormat = "@"

$ws.Range("A5").Value = "Aahil Alwani"
$ws.Range("B5").Value = "4220109168379"
$ws.Range("C5").Value = "2025-05-05T08:18:11.498Z | 2025-05-05T08:19:42.748438Z | 2025-05-05T08:49:27.254322Z"
$ws.Range("D5").Value = "2025-05-05T08:51:19.823867Z | 2025-05-05T08:53:50.148984Z"
$ws.Range("E5").Value = "abbcad36-1f96-4c68-a2a2-ef484605714c"
